$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Drop the disconnected scratch numbers that used to live in L6:N7
#    (L6 held "512 Depth Buffer", L7/M7 held raw numbers, N7 a ratio formula).
# ---------------------------------------------------------------------------
$ws.Range("L6:N7").ClearContents()

# ---------------------------------------------------------------------------
# 2) Insert two rows above the old row 15 ("gaussian_filter..." block).
#    This pushes the old rows 15-17 down to 17-19 and duplicates row 14's
#    formatting onto both freshly inserted rows (Excel's normal Insert
#    behaviour), which is exactly the styling the new spacer row (15) needs.
# ---------------------------------------------------------------------------
$ws.Range("A15:A16").EntireRow.Insert() | Out-Null

# Row 16 is a new sub-header ("benchmark (2048, 2048)") - it should NOT carry
# the copied formatting beyond column B, so clear the extra copied styling.
$ws.Range("C16:J16").Clear() | Out-Null
$ws.Range("B16").ClearFormats() | Out-Null
$ws.Range("B16").Font.Bold = $true

# The moved-down header row (was 15, now 17) switches its B cell style from
# the plain Lucida-Console-left style to Lucida-Console-left+top.
$ws.Range("B17").Font.Name = "Lucida Console"
$ws.Range("B17").HorizontalAlignment = -4131
$ws.Range("B17").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 3) Append the new "Explorer (512, 512)" block (rows 20-22) below the old
#    "full implementation" row (now row 19).
#    Write the new shared-string text in the same order the original commit
#    introduced it (GaussianBlur1, GaussianBlur1Prime, then the two section
#    headers) so the rebuilt shared-string table lines up slot-for-slot.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = "GaussianBlur1"
$ws.Range("C22").Value = "GaussianBlur1Prime"
$ws.Range("B16").Value = "benchmark (2048, 2048)"
$ws.Range("B20").Font.Bold = $true
$ws.Range("B20").Value = "Explorer (512, 512)"

$ws.Range("B21").Font.Name = "Lucida Console"
$ws.Range("B21").HorizontalAlignment = -4131
$ws.Range("B21").Value = 0.011
$ws.Range("F21").Value = 10.556
$ws.Range("H21").Formula = "=F21/`$B`$21"
$ws.Range("H21").NumberFormat = "0.00"

$ws.Range("F22").Value = 7.098
$ws.Range("H22").Formula = "=F22/`$B`$21"
$ws.Range("H22").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 4) Restore the active-cell selection to match the new layout.
# ---------------------------------------------------------------------------
$ws.Range("H24").Select() | Out-Null
